$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.609210848808289
$ws.Range("B1").Value = 1.796569466590881
$ws.Range("C1").Value = 1.83873462677002
$ws.Range("D1").Value = 2.355112075805664
$ws.Range("E1").Value = 3.735209703445435
